# Apply the commit's changes:
#   1. Remove the stray empty "B" cells on rows 6, 7 and 10 of "ODI Batting".
#   2. Add a new worksheet "ODI Batting Extra" (after "ODI Bowling") with
#      per-match batting-extras data (MATCH_CODE, BATTING_POSITION, NUM_4,
#      NUM_6, PERCENT_RUNS_OF_TOTAL, MAN_OF_MATCH).

$wb = $excel.ActiveWorkbook

# --- 1. Clear the stray empty inline-string cells in "ODI Batting" ---------
$odiBatting = $wb.Worksheets.Item("ODI Batting")
$odiBatting.Range("B6").ClearContents()
$odiBatting.Range("B7").ClearContents()
$odiBatting.Range("B10").ClearContents()

# --- 2. Add the new "ODI Batting Extra" worksheet at the end ---------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "ODI Batting Extra"

# Header row
$newSheet.Range("A1").Value = "MATCH_CODE"
$newSheet.Range("B1").Value = "BATTING_POSITION"
$newSheet.Range("C1").Value = "NUM_4"
$newSheet.Range("D1").Value = "NUM_6"
$newSheet.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$newSheet.Range("F1").Value = "MAN_OF_MATCH"

# Re-use the same header style already used on the other sheets ("ODI
# Batting" row 1 is bold, bordered and centred) instead of building a new one.
$odiBatting.Range("A1:F1").Copy()
$newSheet.Range("A1:F1").PasteSpecial(-4122) # xlPasteFormats

# Helper: write $val into $cellRef as TEXT (matching how this workbook
# stores every other "numeric-looking" field, e.g. match codes / counts /
# percentages) without leaving a lingering "Text" number-format style on
# the cell - write to a scratch cell off to the side, then copy/paste
# *values only* so the destination keeps the default style.
$scratch = $newSheet.Range("Z100")
function Set-TextValue($cellRef, $val) {
    $scratch.NumberFormat = "@"
    $scratch.Value = $val
    $scratch.Copy()
    $newSheet.Range($cellRef).PasteSpecial(-4163) # xlPasteValues
}

# Data rows: MATCH_CODE, BATTING_POSITION, NUM_4, NUM_6, PERCENT_RUNS_OF_TOTAL, MAN_OF_MATCH
# BATTING_POSITION is a genuine number; everything else here is text.
$data = @(
    @("4533", $null, $null, $null, $null, "NO"),
    @("4535", 7, "2", "0", "12.24%", "NO"),
    @("4621", 6, "1", "1", "8.77%", "NO"),
    @("4623", $null, $null, $null, $null, "NO"),
    @("4624", 6, $null, $null, $null, "NO"),
    @("4637", 5, $null, $null, $null, "NO"),
    @("4640", 5, "3", "0", "14.97%", "NO"),
    @("4643", $null, $null, $null, $null, "NO"),
    @("4673", $null, $null, $null, $null, "NO"),
    @("4676", 6, "0", "0", "5.48%", "NO")
)

$row = 2
foreach ($rec in $data) {
    Set-TextValue "A$row" $rec[0]
    if ($rec[1] -ne $null) {
        $newSheet.Cells.Item($row, 2).Value = $rec[1]
    }
    if ($rec[2] -ne $null) {
        Set-TextValue "C$row" $rec[2]
    }
    if ($rec[3] -ne $null) {
        Set-TextValue "D$row" $rec[3]
    }
    if ($rec[4] -ne $null) {
        Set-TextValue "E$row" $rec[4]
    }
    Set-TextValue "F$row" $rec[5]
    $row++
}

$scratch.Clear()
